$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-touching the previous last row's phone number makes Excel settle it
# into a real number on save (matches the diff: A60 changes from an
# inline string to a numeric cell once another row is appended after it).
$ws.Cells.Item(60, 1).Value = "76442781"

# Append the new payment row (row 61).
$row = 61

# Write the phone number as text (it's a phone/account id, not a real
# quantity) by forcing a text format before assigning the numeric-looking
# string, then drop back to the default style so no stray formatting is
# left on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "76442781"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = ""
$ws.Cells.Item($row, 3).Value = "Cash"
$ws.Cells.Item($row, 4).Value = "2025-08-18T18:08:21"
$ws.Cells.Item($row, 5).Value = 120
$ws.Cells.Item($row, 6).Value = ""
$ws.Cells.Item($row, 7).Value = 105
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 250
$ws.Cells.Item($row, 10).Value = 15
